$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column C entirely (X_test -> X_train removes the C column values)
$ws.Range("C1:C36").ClearContents()

# Update column D values (Y_test -> Y_train results)
$ws.Range("D1").Value = 0.7845659163987139
$ws.Range("D2").Value = 0.7154340836012861
$ws.Range("D3").Value = 0.9967845659163987
$ws.Range("D4").Value = 0.8697749196141479
$ws.Range("D5").Value = 0.7893890675241158
$ws.Range("D6").Value = 0.6639871382636656
$ws.Range("D7").Value = 0.8164556962025317
$ws.Range("D8").Value = 0.7009493670886076
$ws.Range("D9").Value = 0.9778481012658228
$ws.Range("D10").Value = 0.8987341772151899
$ws.Range("D11").Value = 0.7958860759493671
$ws.Range("D12").Value = 0.685126582278481
$ws.Range("D13").Value = 0.817629179331307
$ws.Range("D14").Value = 0.6975683890577508
$ws.Range("D15").Value = 0.9893617021276596
$ws.Range("D16").Value = 0.9194528875379939
$ws.Range("D17").Value = 0.8024316109422492
$ws.Range("D18").Value = 0.6854103343465046
$ws.Range("D19").Value = 0.7843137254901961
$ws.Range("D20").Value = 0.6802413273001509
$ws.Range("D21").Value = 0.9653092006033183
$ws.Range("D22").Value = 0.8929110105580694
$ws.Range("D23").Value = 0.7873303167420814
$ws.Range("D24").Value = 0.6892911010558069
$ws.Range("D25").Value = 0.7966360856269113
$ws.Range("D26").Value = 0.7064220183486238
$ws.Range("D27").Value = 0.9892966360856269
$ws.Range("D28").Value = 0.8853211009174312
$ws.Range("D29").Value = 0.7920489296636085
$ws.Range("D30").Value = 0.6758409785932722
$ws.Range("D31").Value = 0.810126582278481
$ws.Range("D32").Value = 0.6930379746835443
$ws.Range("D33").Value = 0.9762658227848101
$ws.Range("D34").Value = 0.8718354430379747
$ws.Range("D35").Value = 0.8006329113924051
$ws.Range("D36").Value = 0.6598101265822784

